$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 38360
$ws.Range("D2").Value = 55470947

$ws.Range("C3").Value = 92192
$ws.Range("D3").Value = 135140846

$ws.Range("C4").Value = 31512
$ws.Range("D4").Value = 46668098

$ws.Range("C5").Value = 8835
$ws.Range("D5").Value = 13131356

$ws.Range("C6").Value = 2034
$ws.Range("D6").Value = 3022971

$ws.Range("C12").Value = 41819
$ws.Range("D12").Value = 56729965

$ws.Range("C13").Value = 9810
$ws.Range("D13").Value = 14189991

$ws.Range("C14").Value = 26259
$ws.Range("D14").Value = 38503434

$ws.Range("C15").Value = 8395
$ws.Range("D15").Value = 12457978

$ws.Range("C16").Value = 2181
$ws.Range("D16").Value = 3241435

$ws.Range("C20").Value = 10342
$ws.Range("D20").Value = 13682584

$ws.Range("C21").Value = 13575
$ws.Range("D21").Value = 19592875

$ws.Range("C22").Value = 31987
$ws.Range("D22").Value = 46940717

$ws.Range("C23").Value = 10316
$ws.Range("D23").Value = 15334849

$ws.Range("C24").Value = 2670
$ws.Range("D24").Value = 3969771

$ws.Range("C25").Value = 519
$ws.Range("D25").Value = 772592

$ws.Range("C27").Value = 11826
$ws.Range("D27").Value = 15789672

$ws.Range("C28").Value = 7789
$ws.Range("D28").Value = 11276187

$ws.Range("C29").Value = 22813
$ws.Range("D29").Value = 33484691

$ws.Range("C30").Value = 7888
$ws.Range("D30").Value = 11733891

$ws.Range("C31").Value = 1993
$ws.Range("D31").Value = 2973999

$ws.Range("C32").Value = 371
$ws.Range("D32").Value = 553915

$ws.Range("C34").Value = 8410
$ws.Range("D34").Value = 11109060

$ws.Range("C35").Value = 3310
$ws.Range("D35").Value = 4780653

$ws.Range("C36").Value = 7960
$ws.Range("D36").Value = 11624957

$ws.Range("C37").Value = 3215
$ws.Range("D37").Value = 4765461

$ws.Range("C38").Value = 834
$ws.Range("D38").Value = 1242223

$ws.Range("C41").Value = 2510
$ws.Range("D41").Value = 3393059

$ws.Range("C42").Value = 17512
$ws.Range("D42").Value = 25323011

$ws.Range("C43").Value = 51775
$ws.Range("D43").Value = 75892897

$ws.Range("C44").Value = 19172
$ws.Range("D44").Value = 28475781

$ws.Range("C45").Value = 5673
$ws.Range("D45").Value = 8444760

$ws.Range("C46").Value = 1230
$ws.Range("D46").Value = 1835545

$ws.Range("C47").Value = 65
$ws.Range("D47").Value = 95568

$ws.Range("C50").Value = 16942
$ws.Range("D50").Value = 22531677

$ws.Range("C51").Value = 2093
$ws.Range("D51").Value = 3035678

$ws.Range("C52").Value = 7089
$ws.Range("D52").Value = 10418837

$ws.Range("C53").Value = 2399
$ws.Range("D53").Value = 3582964

$ws.Range("C54").Value = 763
$ws.Range("D54").Value = 1139805

$ws.Range("C55").Value = 189
$ws.Range("D55").Value = 280333

$ws.Range("C57").Value = 7191
$ws.Range("D57").Value = 9887887

$ws.Range("C58").Value = 1077
$ws.Range("D58").Value = 1746006

$ws.Range("C59").Value = 2693
$ws.Range("D59").Value = 4392127

$ws.Range("C60").Value = 1064
$ws.Range("D60").Value = 1739338

$ws.Range("C61").Value = 362
$ws.Range("D61").Value = 594883

$ws.Range("C64").Value = 1579
$ws.Range("D64").Value = 2387992

$ws.Range("C65").Value = 15629
$ws.Range("D65").Value = 22573665

$ws.Range("C66").Value = 45293
$ws.Range("D66").Value = 66271732

$ws.Range("C67").Value = 15875
$ws.Range("D67").Value = 23588566

$ws.Range("C68").Value = 4616
$ws.Range("D68").Value = 6875051

$ws.Range("C69").Value = 944
$ws.Range("D69").Value = 1404168

$ws.Range("C73").Value = 15270
$ws.Range("D73").Value = 20116368

$ws.Range("C74").Value = 53252
$ws.Range("D74").Value = 77498212

$ws.Range("C75").Value = 149911
$ws.Range("D75").Value = 220859737

$ws.Range("C76").Value = 64842
$ws.Range("D76").Value = 96622993

$ws.Range("C77").Value = 20739
$ws.Range("D77").Value = 30988822

$ws.Range("C78").Value = 4938
$ws.Range("D78").Value = 7375901

$ws.Range("C79").Value = 270
$ws.Range("D79").Value = 400170

$ws.Range("C85").Value = 52527
$ws.Range("D85").Value = 71409652

$ws.Range("C86").Value = 4701
$ws.Range("D86").Value = 6813531

$ws.Range("C87").Value = 11760
$ws.Range("D87").Value = 17275888

$ws.Range("C88").Value = 3929
$ws.Range("D88").Value = 5856083

$ws.Range("C89").Value = 1361
$ws.Range("D89").Value = 2033989

$ws.Range("C93").Value = 5495
$ws.Range("D93").Value = 7387423

$ws.Range("C94").Value = 1626
$ws.Range("D94").Value = 2342533

$ws.Range("C95").Value = 5278
$ws.Range("D95").Value = 7774801

$ws.Range("C96").Value = 1962
$ws.Range("D96").Value = 2921726

$ws.Range("C97").Value = 701
$ws.Range("D97").Value = 1050460

$ws.Range("C98").Value = 191
$ws.Range("D98").Value = 287113

$ws.Range("C101").Value = 3631
$ws.Range("D101").Value = 4803940

$ws.Range("C102").Value = 683
$ws.Range("D102").Value = 1109149

$ws.Range("C103").Value = 406
$ws.Range("D103").Value = 668597

$ws.Range("C104").Value = 144
$ws.Range("D104").Value = 232160

$ws.Range("C107").Value = 10963
$ws.Range("D107").Value = 15906568

$ws.Range("C108").Value = 29556
$ws.Range("D108").Value = 43414842

$ws.Range("C109").Value = 9883
$ws.Range("D109").Value = 14695414

$ws.Range("C110").Value = 2721
$ws.Range("D110").Value = 4056580

$ws.Range("C111").Value = 500
$ws.Range("D111").Value = 745046

$ws.Range("C114").Value = 9914
$ws.Range("D114").Value = 13095594

$ws.Range("C115").Value = 30946
$ws.Range("D115").Value = 44619635

$ws.Range("C116").Value = 66929
$ws.Range("D116").Value = 97939943

$ws.Range("C117").Value = 21584
$ws.Range("D117").Value = 32077170

$ws.Range("C118").Value = 6122
$ws.Range("D118").Value = 9121021

$ws.Range("C119").Value = 1146
$ws.Range("D119").Value = 1712600

$ws.Range("C120").Value = 81
$ws.Range("D120").Value = 118920

$ws.Range("C124").Value = 26159
$ws.Range("D124").Value = 34924605

$ws.Range("C125").Value = 36605
$ws.Range("D125").Value = 52824799

$ws.Range("C126").Value = 77787
$ws.Range("D126").Value = 113739192

$ws.Range("C127").Value = 24098
$ws.Range("D127").Value = 35767287

$ws.Range("C128").Value = 6455
$ws.Range("D128").Value = 9593238

$ws.Range("C129").Value = 1260
$ws.Range("D129").Value = 1873311

$ws.Range("C133").Value = 32191
$ws.Range("D133").Value = 42735083

$ws.Range("C134").Value = 13476
$ws.Range("D134").Value = 19510680

$ws.Range("C135").Value = 32714
$ws.Range("D135").Value = 48043938

$ws.Range("C136").Value = 11590
$ws.Range("D136").Value = 17219087

$ws.Range("C137").Value = 2993
$ws.Range("D137").Value = 4460741

$ws.Range("C141").Value = 10931
$ws.Range("D141").Value = 14572825

$ws.Range("C142").Value = 35718
$ws.Range("D142").Value = 51588822

$ws.Range("C143").Value = 82459
$ws.Range("D143").Value = 120807938

$ws.Range("C144").Value = 24643
$ws.Range("D144").Value = 36612763

$ws.Range("C145").Value = 6476
$ws.Range("D145").Value = 9663567

$ws.Range("C146").Value = 1463
$ws.Range("D146").Value = 2176730

$ws.Range("C149").Value = 29593
$ws.Range("D149").Value = 39899942
